# Update scripts with new TPM values for Vtn-Plaur LR-pair sheet
# Adds a new "Resolving-Mac" sending-cluster block (rows 14-17) and
# refreshes numeric columns E:T for every existing row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @('ECs', 'Vtn', 'Plaur', 'ECs', 3, 1, 2.767552, 8.302655999999999, 0.04706493447833917, 0.04706493447833917, 3, 1, 11.128273, 33.384819, 0.1975004092010595, 0.1975004092010595, 30.798074197696, 277.182667779264, 0.009295343818493038, 0.009295343818493038)
$row3 = @('ECs', 'Vtn', 'Plaur', 'FAPs', 3, 1, 2.767552, 8.302655999999999, 0.04706493447833917, 0.04706493447833917, 3, 1, 7.598228333333334, 22.794685, 0.1348505024127659, 0.1348505024127659, 21.02849202037333, 189.25642818336, 0.006346730060427945, 0.006346730060427945)
$row4 = @('ECs', 'Vtn', 'Plaur', 'MuSCs', 3, 1, 2.767552, 8.302655999999999, 0.04706493447833917, 0.04706493447833917, 3, 1, 8.642352, 25.927056, 0.153381217054937, 0.153381217054937, 23.918158562304, 215.263427060736, 0.007218876930898527, 0.007218876930898527)
$row5 = @('ECs', 'Vtn', 'Plaur', 'Resolving-Mac', 3, 1, 2.767552, 8.302655999999999, 0.04706493447833917, 0.04706493447833917, 3, 1, 28.976716, 86.930148, 0.5142678713312377, 0.5142678713312377, 80.19456831923199, 721.7511148730879, 0.02420398366851966, 0.02420398366851966)
$row6 = @('FAPs', 'Vtn', 'Plaur', 'ECs', 3, 1, 21.05317333333333, 63.15952, 0.3580298485789791, 0.3580298485789791, 3, 1, 11.128273, 33.384819, 0.1975004092010595, 0.1975004092010595, 234.2854603696534, 2108.56914332688, 0.07071104160054173, 0.07071104160054173)
$row7 = @('FAPs', 'Vtn', 'Plaur', 'FAPs', 3, 1, 21.05317333333333, 63.15952, 0.3580298485789791, 0.3580298485789791, 3, 1, 7.598228333333334, 22.794685, 0.1348505024127659, 0.1348505024127659, 159.9668181279111, 1439.7013631512, 0.04828050495964183, 0.04828050495964183)
$row8 = @('FAPs', 'Vtn', 'Plaur', 'MuSCs', 3, 1, 21.05317333333333, 63.15952, 0.3580298485789791, 0.3580298485789791, 3, 1, 8.642352, 25.927056, 0.153381217054937, 0.153381217054937, 181.94893466368, 1637.54041197312, 0.05491505391703862, 0.05491505391703861)
$row9 = @('FAPs', 'Vtn', 'Plaur', 'Resolving-Mac', 3, 1, 21.05317333333333, 63.15952, 0.3580298485789791, 0.3580298485789791, 3, 1, 28.976716, 86.930148, 0.5142678713312377, 0.5142678713312377, 610.0518245787733, 5490.46642120896, 0.184123248101757, 0.184123248101757)
$row10 = @('MuSCs', 'Vtn', 'Plaur', 'ECs', 3, 1, 34.97741266666667, 104.932238, 0.5948251867999219, 0.5948251867999219, 3, 1, 11.128273, 33.384819, 0.1975004092010595, 0.1975004092010595, 389.2381969883247, 3503.143772894922, 0.1174782177960812, 0.1174782177960812)
$row11 = @('MuSCs', 'Vtn', 'Plaur', 'FAPs', 3, 1, 34.97741266666667, 104.932238, 0.5948251867999219, 0.5948251867999219, 3, 1, 7.598228333333334, 22.794685, 0.1348505024127659, 0.1348505024127659, 265.766367950559, 2391.89731155503, 0.08021247528773678, 0.08021247528773678)
$row12 = @('MuSCs', 'Vtn', 'Plaur', 'MuSCs', 3, 1, 34.97741266666667, 104.932238, 0.5948251867999219, 0.5948251867999219, 3, 1, 8.642352, 25.927056, 0.153381217054937, 0.153381217054937, 302.2871123145921, 2720.584010831328, 0.09123501108630226, 0.09123501108630225)
$row13 = @('MuSCs', 'Vtn', 'Plaur', 'Resolving-Mac', 3, 1, 34.97741266666667, 104.932238, 0.5948251867999219, 0.5948251867999219, 3, 1, 28.976716, 86.930148, 0.5142678713312377, 0.5142678713312377, 1013.530553256803, 9121.774979311225, 0.3058994826298017, 0.3058994826298017)
$row14 = @('Resolving-Mac', 'Vtn', 'Plaur', 'ECs', 1, 0.3333333333333333, 0.004706, 0.014118, 0.00008003014275976175, 0.00008003014275976175, 3, 1, 11.128273, 33.384819, 0.1975004092010595, 0.1975004092010595, 0.052369652738, 0.471326874642, 0.00001580598594347215, 0.00001580598594347215)
$row15 = @('Resolving-Mac', 'Vtn', 'Plaur', 'FAPs', 1, 0.3333333333333333, 0.004706, 0.014118, 0.00008003014275976175, 0.00008003014275976175, 3, 1, 7.598228333333334, 22.794685, 0.1348505024127659, 0.1348505024127659, 0.03575726253666667, 0.32181536283, 0.00001079210495931925, 0.00001079210495931925)
$row16 = @('Resolving-Mac', 'Vtn', 'Plaur', 'MuSCs', 1, 0.3333333333333333, 0.004706, 0.014118, 0.00008003014275976175, 0.00008003014275976175, 3, 1, 8.642352, 25.927056, 0.153381217054937, 0.153381217054937, 0.040670908512, 0.366038176608, 0.00001227512069757261, 0.00001227512069757261)
$row17 = @('Resolving-Mac', 'Vtn', 'Plaur', 'Resolving-Mac', 1, 0.3333333333333333, 0.004706, 0.014118, 0.00008003014275976175, 0.00008003014275976175, 3, 1, 28.976716, 86.930148, 0.5142678713312377, 0.5142678713312377, 0.136364425496, 1.227279829464, 0.00004115693115939775, 0.00004115693115939775)

$rowData = @{
  2 = $row2
  3 = $row3
  4 = $row4
  5 = $row5
  6 = $row6
  7 = $row7
  8 = $row8
  9 = $row9
  10 = $row10
  11 = $row11
  12 = $row12
  13 = $row13
  14 = $row14
  15 = $row15
  16 = $row16
  17 = $row17
}

foreach ($r in ($rowData.Keys | Sort-Object)) {
  $vals = $rowData[$r]
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item([int]$r, $col).Value = $vals[$i]
  }
}
